$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input")

$ws.Range("D2").Value = "Network based on OSM online"
$ws.Range("E2").Value = "zuidholland_4326"
$ws.Range("F2").Value = ""
$ws.Range("L2").Value = "drive"
$ws.Range("M2").Value = "motorway"

$ws.Activate()
$ws.Range("M7").Select()
